$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 15 (Part P-14): Claw Plate, 1/4'' thick, no bending
$ws.Range("C15").Value = "1/4''"
$ws.Range("D15").Value = "no"
$ws.Range("F15").Value = "Claw Plate"

# Row 16 (Part P-15): Claw Slide Plate, 1/8'' thick, "no " bending
$ws.Range("C16").Value = "1/8''"
$ws.Range("D16").Value = "no "
$ws.Range("F16").Value = "Claw Slide Plate"

# Row 17 (Part P-16): Claw Support, 1/8'' thick, yes bending
$ws.Range("C17").Value = "1/8''"
$ws.Range("D17").Value = "yes"
$ws.Range("F17").Value = "Claw Support "

# Row 18 (Part P-17): Claw Support Other Side, 1/8'' thick, yes bending
$ws.Range("C18").Value = "1/8''"
$ws.Range("D18").Value = "yes"
$ws.Range("F18").Value = "Claw Support Other Side"

# Leave selection on D14, matching the final cursor position recorded in the file
$ws.Range("D14").Select()
